# Registration testcase's modification with validation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("advertiserRegistration")

# --- Row 2: "Select Client Type *" / "Select Country *" label cells -------
$fLabelColor = 2236962  # RGB(34,34,34) == FF222222

$f2 = $ws.Cells.Item(2, 6)
$f2.Value = "Select Client Type *"
$f2.Font.Color = $fLabelColor
$f2.Borders.Item(7).LineStyle = -4142
$f2.Borders.Item(8).LineStyle = -4142
$f2.Borders.Item(9).LineStyle = -4142
$f2.Borders.Item(10).LineStyle = -4142
$f2.VerticalAlignment = -4160

$q2 = $ws.Cells.Item(2, 17)
$q2.Value = "Select Country *"
$q2.Font.Color = $fLabelColor
$q2.Borders.Item(7).LineStyle = -4142
$q2.Borders.Item(8).LineStyle = -4142
$q2.Borders.Item(9).LineStyle = -4142
$q2.Borders.Item(10).LineStyle = -4142

# --- Action column (U) wording fixes: cancle/Submit -> cancel/submit ------
$ws.Cells.Item(2, 21).Value = "cancel"
$ws.Cells.Item(3, 21).Value = "cancel"
$ws.Cells.Item(4, 21).Value = "submit"
$ws.Cells.Item(5, 21).Value = "cancel"
$ws.Cells.Item(6, 21).Value = ""

# --- Expected Message column (V) content updates --------------------------
$ws.Cells.Item(2, 22).Value = "Client Type can not be blank, Organization Name cant be blank, First Name can't be blank, Last Name cant be blank, Email Id can't be blank,Password cant be blank, Confirm Password can't be blank, Contact Number can't be blank, Website can't be blank, Address can't be blank, Country Name can't be blank, State Name can't be blank, City Name can't be blank, Pincode Number can't be blank,"
$ws.Cells.Item(2, 22).WrapText = $true

$ws.Cells.Item(3, 22).Value = "Email id is invalid, Confirm Password is not same as Password, Website is invalid,"
$ws.Cells.Item(5, 22).Value = " Password length should be Greater than 8, Confirm Password is not same as Password, Contact Number Invalid, Website is invalid,"

# --- Header "Expected Message" (V1): center + wrap ------------------------
$v1 = $ws.Cells.Item(1, 22)
$v1.HorizontalAlignment = -4108
$v1.VerticalAlignment = -4160
$v1.WrapText = $true

# --- Column V: default wrap-text style + row heights -----------------------
$ws.Columns("V").WrapText = $true

$ws.Rows(2).RowHeight = 60
$ws.Rows(3).RowHeight = 15
$ws.Rows(4).RowHeight = 30
$ws.Rows(5).RowHeight = 30

# --- View state (best effort) ----------------------------------------------
$ws.Range("A1:V6").Select()
